$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("README")

# Clear the previous README body (content + formatting) so it can be rebuilt with the new layout
$ws.Range("A1:M24").Clear() | Out-Null

$ws.Range("A1").Value = 'SCP Sample Annotation Wizard'
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 14

$ws.Range("A2").Value = 'Generates sample annotation .csv files to be passed to the colData argument in the readSCP() function from the scp R package (https://uclouvain-cbio.github.io/scp/index.html)'

$ws.Range("A4").Value = 'User Guidance'
$ws.Range("A4").Font.Bold = $true
$ws.Range("A4").Font.Size = 14

$ws.Range("A5").Value = 'Page'
$ws.Range("A5").Font.Bold = $true
$ws.Range("A5").Font.Size = 11

$ws.Range("B5").Value = 'Form'
$ws.Range("B5").Font.Bold = $true
$ws.Range("B5").Font.Size = 11

$ws.Range("C5").Value = 'Guidance'
$ws.Range("C5").Font.Bold = $true
$ws.Range("C5").Font.Size = 11

$ws.Range("A6").Value = 'Import Page'

$ws.Range("B6").Value = 'Technology used'

$ws.Range("C6").Value = 'Select the technology used to generate the data'

$ws.Range("B7").Value = 'Date File Extension'

$ws.Range("C7").Value = 'Select the type of data used'

$ws.Range("B8").Value = 'Data files'

$ws.Range("C8").Value = 'Import either folder containing ".raw" files or CSV formatted like ''Raw File Import - Template'''

$ws.Range("B9").Value = 'Labels file / droplet location file'

$ws.Range("C9").Value = 'Import .fld'

$ws.Range("B10").Value = 'Pickup file'

$ws.Range("C10").Value = 'Import .log file or .log converted to csv'

$ws.Range("B11").Value = 'Cell files'

$ws.Range("C11").Value = 'Import as many cell files as desired in .xls format, or converted to .csv'

$ws.Range("B12").Value = 'Additional cellenONE annotation files'

$ws.Range("C12").Value = 'Import additional cellenONE annotation files in .fld format to include as an output column'

$ws.Range("B13").Value = 'Task name'

$ws.Range("C13").Value = 'Provide folder name within "documents/ScpSampleAnnotationWizardOutput/" to send output to'

$ws.Range("M13").Style = "Hyperlink"

$ws.Range("A15").Value = 'Metadata Page'

$ws.Range("B15").Value = 'Regex to extract row/ column'

$ws.Range("C15").Value = 'Edit the regex that is used to extract row and column names from raw file names if the default has failed to do so'

$ws.Range("B16").Value = 'Metadata to Include'

$ws.Range("C16").Value = 'Tick columns to include in output file'

$ws.Range("B17").Value = 'Handle Cell File Column Mismatches'

$ws.Range("C17").Value = 'Decide how to handle columns that are present in one cell file, but not others'

$ws.Range("A19").Value = 'Other Options Page'

$ws.Range("B19").Value = 'Regex to extract well'

$ws.Range("C19").Value = 'Enter regex to correctly extract well to merge cell files and data files'

$ws.Range("B20").Value = 'Name cell populations'

$ws.Range("C20").Value = 'Choose how cell files should be named'

$ws.Range("B21").Value = 'Name missing cell data values'

$ws.Range("C21").Value = 'Enter a value to fill any missing row'

$ws.Range("B22").Value = 'Add extra rows per raw file'

$ws.Range("C22").Value = 'Any extra rows to add per raw file'

$ws.Range("B23").Value = 'Well to TMT mapping CSV'

$ws.Range("C23").Value = 'Choose the mapping of well to TMT file, if default selected ''Well to TMT mapping - Default'' is used, else ''Well to TMT mapping - Template'' can be populated and uploaded as a csv'

$ws.Range("B24").Value = 'Well to Label mapping CSV (Only if "Label-based" selected as technology)'

$ws.Range("C24").Value = 'Choose the mapping of well to Label, if ''Default'' selected ''Well to Label mapping - Default'' is used, else ''Well to Label mapping - Template'' can be populated and uploaded as a csv'

$ws.Range("B25").Value = 'Pickup type'

$ws.Range("C25").Value = 'Select whether data was generated from single pickup or dual pickup. If dual, then include the offset of the X position from the first position to the second'
$ws.Range("C25").WrapText = $true

$ws.Range("B26").Value = 'Invert numbering'

$ws.Range("C26").Value = 'Select whether the numbering of column or row should be inverted e.g. "3, 2, 1" instead of "1, 2, 3". And if so, enter the regex to correctly select the column or row.'

$ws.Range("A28").Value = 'Note'
$ws.Range("A28").Font.Bold = $true
$ws.Range("A28").Font.Size = 11

$ws.Range("A29").Value = 'Editing the name of the application will lead to errors'

# Restore the reported selection/active cell for the README sheet
$ws.Range("B23").Select() | Out-Null